$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
  "AB2" = 6.6
  "AF2" = 21
  "F2" = 2.04
  "H2" = 4.4
  "I2" = 5.1
  "J2" = 3
  "K2" = 3.35
  "L2" = 1.62
  "N2" = 2.46
  "O2" = 1.59
  "P2" = 1.47
  "Q2" = 2.84
  "R2" = 1.17
  "S2" = 5.7
  "T2" = 2.22
  "U2" = 1.66
  "V2" = 1.26
  "W2" = 1.83
  "AB3" = 19.5
  "AC3" = 7.4
  "AE3" = 980
  "F3" = 2.94
  "G3" = 3.25
  "H3" = 2.88
  "I3" = 3.15
  "J3" = 2.78
  "L3" = 1.7
  "M3" = 1.16
  "N3" = 2.26
  "O3" = 1.69
  "Q3" = 3.05
  "S3" = 6.8
  "T3" = 2.3
  "U3" = 1.65
  "V3" = 1.47
  "W3" = 1.45
  "X3" = 14
  "AC4" = 26
  "AD4" = 300
  "AE4" = 340
  "AF4" = 10.5
  "AH4" = 48
  "AJ4" = 10
  "AK4" = 13.5
  "AL4" = 34
  "AN4" = 2.72
  "I4" = 26
  "J4" = 9.6
  "L4" = 1.19
  "N4" = 9.6
  "P4" = 3.9
  "Q4" = 1.29
  "R4" = 2.16
  "S4" = 1.78
  "T4" = 1.95
  "U4" = 1.87
  "W4" = 6.6
  "X4" = 990
  "Y4" = 1000
  "Z4" = 260
  "AB5" = 6.2
  "AC5" = 8
  "AE5" = 130
  "AH5" = 29
  "AK5" = 29
  "AL5" = 70
  "AN5" = 26
  "F5" = 1.86
  "K5" = 3.5
  "L5" = 1.6
  "M5" = 1.13
  "N5" = 2.5
  "O5" = 1.56
  "P5" = 1.49
  "Q5" = 2.74
  "R5" = 1.17
  "S5" = 5.7
  "T5" = 2.3
  "U5" = 1.64
  "X5" = 8.800000000000001
  "Y5" = 13.5
  "Z5" = 50
  "AB6" = 8.199999999999999
  "AD6" = 14.5
  "AE6" = 48
  "AG6" = 12.5
  "AH6" = 22
  "AI6" = 75
  "AJ6" = 38
  "AK6" = 34
  "AL6" = 60
  "AM6" = 160
  "AN6" = 36
  "AO6" = 70
  "L6" = 1.57
  "M6" = 1.12
  "N6" = 2.86
  "O6" = 1.52
  "P6" = 1.61
  "Q6" = 2.62
  "R6" = 1.22
  "S6" = 5.3
  "T6" = 2.08
  "U6" = 1.9
  "X6" = 8.800000000000001
  "Y6" = 9.6
  "AB7" = 16.5
  "AC7" = 13.5
  "AD7" = 23
  "AF7" = 15
  "AI7" = 46
  "AL7" = 24
  "AM7" = 60
  "AN7" = 5.4
  "AO7" = 46
  "F7" = 1.58
  "H7" = 5.2
  "L7" = 1.23
  "N7" = 7.6
  "O7" = 1.13
  "P7" = 3.25
  "Q7" = 1.39
  "R7" = 1.9
  "S7" = 2
  "T7" = 1.52
  "U7" = 2.66
  "X7" = 38
  "Y7" = 40
  "AA8" = 900
  "AB8" = 34
  "AC8" = 8.800000000000001
  "AE8" = 38
  "AF8" = 980
  "AH8" = 38
  "AN8" = 1000
  "AO8" = 50
  "F8" = 4.6
  "G8" = 5.1
  "H8" = 1.88
  "I8" = 1.94
  "J8" = 3.65
  "K8" = 3.9
  "L8" = 1.42
  "N8" = 3.75
  "P8" = 1.98
  "Q8" = 1.96
  "S8" = 3.45
  "T8" = 1.81
  "V8" = 2.06
  "W8" = 1.25
  "Y8" = 9
  "Z8" = 12
  "AC9" = 9
  "AK9" = 24
  "F9" = 1.65
  "I9" = 8.199999999999999
  "L9" = 1.56
  "N9" = 2.72
  "Q9" = 2.52
  "R9" = 1.2
  "U9" = 1.61
  "Y9" = 18
  "AA10" = 440
  "AB10" = 6.2
  "AC10" = 10
  "AD10" = 38
  "AE10" = 200
  "AF10" = 7.4
  "AG10" = 10
  "AI10" = 180
  "AJ10" = 12
  "AK10" = 18
  "AL10" = 55
  "AM10" = 270
  "AN10" = 9.4
  "F10" = 1.46
  "G10" = 1.47
  "H10" = 9.4
  "I10" = 9.6
  "J10" = 4.6
  "K10" = 4.7
  "L10" = 1.47
  "N10" = 3.35
  "O10" = 1.4
  "P10" = 1.81
  "Q10" = 2.2
  "R10" = 1.3
  "S10" = 4.2
  "T10" = 2.42
  "U10" = 1.66
  "V10" = 1.11
  "W10" = 3.1
  "X10" = 13
  "Y10" = 24
  "Z10" = 85
  "AA11" = 200
  "AB11" = 7.8
  "AC11" = 9.4
  "AE11" = 80
  "AF11" = 8.800000000000001
  "AH11" = 22
  "AI11" = 80
  "AJ11" = 15.5
  "AK11" = 15.5
  "AL11" = 32
  "AM11" = 110
  "AN11" = 9.4
  "AO11" = 120
  "F11" = 1.59
  "G11" = 1.6
  "H11" = 6.4
  "I11" = 6.6
  "L11" = 1.42
  "M11" = 1.06
  "N11" = 4
  "O11" = 1.31
  "P11" = 2.04
  "Q11" = 1.94
  "R11" = 1.41
  "S11" = 3.35
  "T11" = 1.96
  "U11" = 2
  "W11" = 2.66
  "X11" = 17.5
  "Y11" = 19.5
  "Z11" = 55
  "AA12" = 1000
  "AB12" = 7.8
  "AC12" = 7.6
  "AD12" = 17.5
  "AF12" = 13
  "AG12" = 11.5
  "AJ12" = 29
  "AK12" = 28
  "AL12" = 50
  "AM12" = 580
  "AN12" = 28
  "AO12" = 1000
  "F12" = 2.14
  "G12" = 2.26
  "H12" = 3.95
  "I12" = 4.4
  "J12" = 3.15
  "K12" = 3.45
  "L12" = 1.51
  "N12" = 3
  "O12" = 1.44
  "P12" = 1.66
  "Q12" = 2.36
  "T12" = 2
  "U12" = 1.89
  "V12" = 1.3
  "W12" = 1.79
  "X12" = 10.5
  "Y12" = 12.5
  "Z12" = 29
}

foreach ($addr in $changes.Keys) {
  $ws.Range($addr).Value = $changes[$addr]
}
